# Trade #63 closed at 2026-02-17 21:11:19 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh headline stats now that another trade has closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.11
$summary.Range("B4").Value = 0.91
$summary.Range("B5").Value = 0.2
$summary.Range("B6").Value = 91
$summary.Range("B8").Value = 35
$summary.Range("B9").Value = 48.35

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.11
$status.Range("D5").Value = 58
$status.Range("E5").Value = 0.8
$status.Range("F5").Value = 1.11
$status.Range("G5").Value = 51.72

# ---------------------------------------------------------------------------
# Sheet "All Trades": trade #91 (row 92) flips from OPEN to CLOSED via an
# early exit, and a brand-new trade #124 (row 125) is appended as OPEN.
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G92").Value = 0.09
$allTrades.Range("H92").Value = "CLOSED"
$allTrades.Range("I92").Value = -30.7692
$allTrades.Range("J92").Value = -0.04
$allTrades.Range("K92").Value = 101.11
$allTrades.Range("L92").Value = "early_exit"
$allTrades.Range("M92").Value = 0.14

# New row 125 - keep the Date column ("2026-02-17") stored as literal text,
# not auto-converted to a date serial number, matching the rest of the sheet.
$dateCell = $allTrades.Range("B125")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-02-17"
$dateCell.Style = "Normal"

$allTrades.Range("A125").Value = 124
$allTrades.Range("C125").Value = "21:11:13"
$allTrades.Range("D125").Value = "MarketMaking"
$allTrades.Range("E125").Value = "DOWN"
$allTrades.Range("F125").Value = 0.13
$allTrades.Range("G125").Value = ""
$allTrades.Range("H125").Value = "OPEN"
$allTrades.Range("I125").Value = 0
$allTrades.Range("J125").Value = 0
$allTrades.Range("K125").Value = 101.1546450978375
$allTrades.Range("L125").Value = ""
$allTrades.Range("M125").Value = 0
$allTrades.Range("N125").Value = 0
$allTrades.Range("O125").Value = 0
$allTrades.Range("P125").Value = 0.6
$allTrades.Range("Q125").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Sheet "MarketMaking": the same trade #91 lives at row 59 here (per-strategy
# numbering), and the new trade #124 is appended at row 92.
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Range("G59").Value = 0.09
$marketMaking.Range("H59").Value = "CLOSED"
$marketMaking.Range("I59").Value = -30.7692
$marketMaking.Range("J59").Value = -0.04
$marketMaking.Range("K59").Value = 101.11
$marketMaking.Range("P59").Value = "early_exit"
$marketMaking.Range("Q59").Value = 0.14

$dateCell2 = $marketMaking.Range("B92")
$dateCell2.NumberFormat = "@"
$dateCell2.Value = "2026-02-17"
$dateCell2.Style = "Normal"

$marketMaking.Range("A92").Value = 124
$marketMaking.Range("C92").Value = "21:11:13"
$marketMaking.Range("D92").Value = "MarketMaking"
$marketMaking.Range("E92").Value = "DOWN"
$marketMaking.Range("F92").Value = 0.13
$marketMaking.Range("G92").Value = ""
$marketMaking.Range("H92").Value = "OPEN"
$marketMaking.Range("I92").Value = 0
$marketMaking.Range("J92").Value = 0
$marketMaking.Range("K92").Value = 101.1546450978375
$marketMaking.Range("L92").Value = 0
$marketMaking.Range("M92").Value = 0
$marketMaking.Range("N92").Value = 0.6
$marketMaking.Range("O92").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("P92").Value = ""
$marketMaking.Range("Q92").Value = 0
